$d = $word.ActiveDocument

# 1. "CSCI UA.0060 Spring 2024" -> "CSCI UA.0060 Fall 2024"
$d.Content.Find.Execute("Spring 2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fall 2024", 2)

# 2. Requirements paragraph: describe LucidChart workflow instead of "any tool of your choice"
$d.Content.Find.Execute( `
    "The diagram can be created using Visio, Excel or any other tool of your choice, but it must be saved/exported as a PDF", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "The diagram should be created using LucidChart and then saved as a PDF", 2)
